$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = 0.1586660023084538
$ws.Range("C2").Value = 1.986674183974745
$ws.Range("D2").Value = 17.01968672753467
$ws.Range("E2").Value = 4.125492301233233
$ws.Range("F2").Value = 4.172410827772781
$ws.Range("G2").Value = 42

# Row 3 (Q0)
$ws.Range("B3").Value = 0.2379729132671037
$ws.Range("C3").Value = 1.913853995852431
$ws.Range("D3").Value = 14.70470842293039
$ws.Range("E3").Value = 3.834671879435109
$ws.Range("F3").Value = 3.841122596751489
$ws.Range("G3").Value = 139

# Row 4 (Q1)
$ws.Range("B4").Value = 0.1721244478012908
$ws.Range("C4").Value = 1.272346930490598
$ws.Range("D4").Value = 5.520684595974761
$ws.Range("E4").Value = 2.34961371207583
$ws.Range("F4").Value = 2.3602199723419
$ws.Range("G4").Value = 70
